$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 1621
$ws.Range("E2").Value = -38
$ws.Range("F2").Value = -38
$ws.Range("G2").Value = -56
$ws.Range("H2").Value = -60
$ws.Range("I2").Value = -62
$ws.Range("J2").Value = 2
$ws.Range("K2").Value = 1088
$ws.Range("L2").Value = 629
$ws.Range("M2").Value = 460
$ws.Range("N2").Value = 447
$ws.Range("O2").Value = 13
$ws.Range("P2").Value = 104
$ws.Range("Q2").Value = 37
$ws.Range("R2").Value = -16
$ws.Range("S2").Value = 7
$ws.Range("T2").Value = 30
$ws.Range("U2").Value = 6
$ws.Range("V2").Value = 322
$ws.Range("W2").Value = -2.37
$ws.Range("X2").Value = -3.71
$ws.Range("Y2").Value = -12.76
$ws.Range("Z2").Value = -5.16
$ws.Range("AA2").Value = 136.67
$ws.Range("AB2").Value = 292.72
$ws.Range("AC2").Value = -598
$ws.Range("AD2").Value = -6.23
$ws.Range("AE2").Value = 4356
$ws.Range("AF2").Value = 0.85
$ws.Range("AG2").Value = 50
$ws.Range("AH2").Value = 1.34
$ws.Range("AI2").Value = -6.67
$ws.Range("AJ2").Value = 10395000

# Row 3
$ws.Range("D3").Value = 1588
$ws.Range("E3").Value = 54
$ws.Range("F3").Value = 54
$ws.Range("G3").Value = 34
$ws.Range("H3").Value = 29
$ws.Range("I3").Value = 27
$ws.Range("J3").Value = 2
$ws.Range("K3").Value = 1129
$ws.Range("L3").Value = 657
$ws.Range("M3").Value = 473
$ws.Range("N3").Value = 465
$ws.Range("O3").Value = 8
$ws.Range("P3").Value = 104
$ws.Range("Q3").Value = 138
$ws.Range("R3").Value = -77
$ws.Range("S3").Value = -77
$ws.Range("T3").Value = 38
$ws.Range("U3").Value = 99
$ws.Range("V3").Value = 293
$ws.Range("W3").Value = 3.4
$ws.Range("X3").Value = 1.8
$ws.Range("Y3").Value = 5.85
$ws.Range("Z3").Value = 2.58
$ws.Range("AA3").Value = 139.02
$ws.Range("AB3").Value = 301.17
$ws.Range("AC3").Value = 257
$ws.Range("AD3").Value = 45.78
$ws.Range("AE3").Value = 4529
$ws.Range("AF3").Value = 2.59
$ws.Range("AG3").Value = 50
$ws.Range("AH3").Value = 0.43
$ws.Range("AI3").Value = 19.23
$ws.Range("AJ3").Value = 10395000

# Row 4
$ws.Range("D4").Value = 1725
$ws.Range("E4").Value = 95
$ws.Range("F4").Value = 95
$ws.Range("G4").Value = 72
$ws.Range("H4").Value = 57
$ws.Range("I4").Value = 56
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 1217
$ws.Range("L4").Value = 683
$ws.Range("M4").Value = 534
$ws.Range("N4").Value = 526
$ws.Range("O4").Value = 8
$ws.Range("P4").Value = 104
$ws.Range("Q4").Value = 127
$ws.Range("R4").Value = -53
$ws.Range("S4").Value = -53
$ws.Range("T4").Value = 34
$ws.Range("U4").Value = 93
$ws.Range("V4").Value = 311
$ws.Range("W4").Value = 5.49
$ws.Range("X4").Value = 3.29
$ws.Range("Y4").Value = 11.31
$ws.Range("Z4").Value = 4.84
$ws.Range("AA4").Value = 128.02
$ws.Range("AB4").Value = 350.94
$ws.Range("AC4").Value = 539
$ws.Range("AD4").Value = 19.86
$ws.Range("AE4").Value = 5123
$ws.Range("AF4").Value = 2.09
$ws.Range("AG4").Value = 50
$ws.Range("AH4").Value = 0.47
$ws.Range("AI4").Value = 9.16
$ws.Range("AJ4").Value = 10395000

# Row 5
$ws.Range("D5").Value = 1993
$ws.Range("E5").Value = 220
$ws.Range("F5").Value = 220
$ws.Range("G5").Value = 201
$ws.Range("H5").Value = 155
$ws.Range("I5").Value = 155
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 1377
$ws.Range("L5").Value = 687
$ws.Range("M5").Value = 690
$ws.Range("N5").Value = 681
$ws.Range("O5").Value = 9
$ws.Range("P5").Value = 104
$ws.Range("Q5").Value = 207
$ws.Range("R5").Value = -112
$ws.Range("S5").Value = -60
$ws.Range("T5").Value = 107
$ws.Range("U5").Value = 101
$ws.Range("V5").Value = 205
$ws.Range("W5").Value = 11.06
$ws.Range("X5").Value = 7.8
$ws.Range("Y5").Value = 25.71
$ws.Range("Z5").Value = 11.99
$ws.Range("AA5").Value = 99.58
$ws.Range("AB5").Value = 495.09
$ws.Range("AC5").Value = 1493
$ws.Range("AD5").Value = 28.84
$ws.Range("AE5").Value = 6641
$ws.Range("AF5").Value = 6.48
$ws.Range("AG5").Value = 150
$ws.Range("AH5").Value = 0.35
$ws.Range("AI5").Value = 9.92
$ws.Range("AJ5").Value = 10395000

# Row 6
$ws.Range("D6").Value = 2744
$ws.Range("E6").Value = 813
$ws.Range("F6").Value = 813
$ws.Range("G6").Value = 812
$ws.Range("H6").Value = 619
$ws.Range("I6").Value = 618
$ws.Range("K6").Value = 2200
$ws.Range("L6").Value = 927
$ws.Range("M6").Value = 1274
$ws.Range("N6").Value = 1264
$ws.Range("P6").Value = 104
$ws.Range("Q6").Value = 609
$ws.Range("R6").Value = -431
$ws.Range("S6").Value = -35
$ws.Range("T6").Value = 277
$ws.Range("U6").Value = 332
$ws.Range("V6").Value = 201
$ws.Range("W6").Value = 29.64
$ws.Range("X6").Value = 22.56
$ws.Range("Y6").Value = 63.52
$ws.Range("Z6").Value = 34.62
$ws.Range("AA6").Value = 72.76
$ws.Range("AB6").Value = 1061.41
$ws.Range("AC6").Value = 5945
$ws.Range("AD6").Value = 8.95
$ws.Range("AE6").Value = 12321
$ws.Range("AF6").Value = 4.32
$ws.Range("AG6").Value = 300
$ws.Range("AH6").Value = 0.56
$ws.Range("AI6").Value = 4.98
$ws.Range("AJ6").Value = 10395000

# Row 7
$ws.Range("D7").Value = 2578
$ws.Range("E7").Value = 439
$ws.Range("G7").Value = 451
$ws.Range("H7").Value = 358
$ws.Range("I7").Value = 357
$ws.Range("K7").Value = 2363
$ws.Range("L7").Value = 750
$ws.Range("M7").Value = 1610
$ws.Range("N7").Value = 1597
$ws.Range("P7").Value = 100
$ws.Range("Q7").Value = 507
$ws.Range("R7").Value = -283
$ws.Range("S7").Value = -40
$ws.Range("T7").Value = 535
$ws.Range("U7").Value = -65
$ws.Range("W7").Value = 17.02
$ws.Range("X7").Value = 13.9
$ws.Range("Y7").Value = 24.93
$ws.Range("Z7").Value = 15.7
$ws.Range("AA7").Value = 46.58
$ws.Range("AC7").Value = 3431
$ws.Range("AD7").Value = 15.88
$ws.Range("AE7").Value = 15561
$ws.Range("AF7").Value = 3.5
$ws.Range("AG7").Value = 200
$ws.Range("AH7").Value = 0.37
$ws.Range("AI7").Value = 5.83

# Row 8
$ws.Range("D8").Value = 2843
$ws.Range("E8").Value = 551
$ws.Range("G8").Value = 560
$ws.Range("H8").Value = 433
$ws.Range("I8").Value = 427
$ws.Range("K8").Value = 2823
$ws.Range("L8").Value = 803
$ws.Range("M8").Value = 2023
$ws.Range("N8").Value = 2007
$ws.Range("P8").Value = 100
$ws.Range("Q8").Value = 627
$ws.Range("R8").Value = -143
$ws.Range("S8").Value = -47
$ws.Range("T8").Value = 200
$ws.Range("U8").Value = 460
$ws.Range("W8").Value = 19.38
$ws.Range("X8").Value = 15.24
$ws.Range("Y8").Value = 23.68
$ws.Range("Z8").Value = 16.71
$ws.Range("AA8").Value = 39.7
$ws.Range("AC8").Value = 4105
$ws.Range("AD8").Value = 13.28
$ws.Range("AE8").Value = 19556
$ws.Range("AF8").Value = 2.79
$ws.Range("AG8").Value = 233
$ws.Range("AH8").Value = 0.43
$ws.Range("AI8").Value = 5.68

# Row 9
$ws.Range("D9").Value = 3263
$ws.Range("E9").Value = 623
$ws.Range("G9").Value = 637
$ws.Range("H9").Value = 490
$ws.Range("I9").Value = 490
$ws.Range("K9").Value = 3333
$ws.Range("L9").Value = 847
$ws.Range("M9").Value = 2487
$ws.Range("N9").Value = 2473
$ws.Range("P9").Value = 100
$ws.Range("Q9").Value = 660
$ws.Range("R9").Value = -200
$ws.Range("S9").Value = -23
$ws.Range("T9").Value = 300
$ws.Range("U9").Value = 450
$ws.Range("W9").Value = 19.1
$ws.Range("X9").Value = 15.02
$ws.Range("Y9").Value = 21.88
$ws.Range("Z9").Value = 15.92
$ws.Range("AA9").Value = 34.05
$ws.Range("AC9").Value = 4714
$ws.Range("AD9").Value = 11.56
$ws.Range("AE9").Value = 24104
$ws.Range("AF9").Value = 2.26
$ws.Range("AG9").Value = 267
$ws.Range("AH9").Value = 0.49
$ws.Range("AI9").Value = 5.66
